# UPD data from 2020-05-10
# Adds two new dates (2020-05-08, 2020-05-09) of "si" data for each of the
# three states (Moscow, Sevastopol, Saint Petersburg), inserting new rows
# right after each state's existing block in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText($cell, $text) {
    # Force the date-like string to be stored as text (not auto-converted
    # to a date serial number), then restore the cell's default style so
    # no extra formatting is applied to the cell.
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = $ws.Cells.Item(1,1).Style
}

# --- Insert 2 rows after Moscow block (rows 2-76) -> new rows 77,78 ---
$ws.Rows.Item(77).Resize(2).Insert()
$ws.Cells.Item(77,1).Value2 = "Москва"
Set-DateText $ws.Cells.Item(77,2) "2020-05-08"
$ws.Cells.Item(77,3).Value2 = 2.8
$ws.Cells.Item(78,1).Value2 = "Москва"
Set-DateText $ws.Cells.Item(78,2) "2020-05-09"
$ws.Cells.Item(78,3).Value2 = 3.6

# --- Insert 2 rows after Sevastopol block (now rows 79-153) -> new rows 154,155 ---
$ws.Rows.Item(154).Resize(2).Insert()
$ws.Cells.Item(154,1).Value2 = "Севастополь"
Set-DateText $ws.Cells.Item(154,2) "2020-05-08"
$ws.Cells.Item(154,3).Value2 = 2.8
$ws.Cells.Item(155,1).Value2 = "Севастополь"
Set-DateText $ws.Cells.Item(155,2) "2020-05-09"
$ws.Cells.Item(155,3).Value2 = 3.3

# --- Append 2 rows after Saint Petersburg block (now rows 156-230) -> new rows 231,232 ---
$ws.Cells.Item(231,1).Value2 = "Санкт-Петербург"
Set-DateText $ws.Cells.Item(231,2) "2020-05-08"
$ws.Cells.Item(231,3).Value2 = 2.2
$ws.Cells.Item(232,1).Value2 = "Санкт-Петербург"
Set-DateText $ws.Cells.Item(232,2) "2020-05-09"
$ws.Cells.Item(232,3).Value2 = 3.1
